$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "model_8_1_21"
$ws.Range("B2").Value2 = 0.2967959353427138
$ws.Range("C2").Value2 = -5.092416589398046
$ws.Range("D2").Value2 = -2.262551119888867
$ws.Range("E2").Value2 = -3.236308560062831
$ws.Range("F2").Value2 = 0.7782397270202637
$ws.Range("G2").Value2 = 2.761625051498413
$ws.Range("H2").Value2 = 2.323813199996948
$ws.Range("I2").Value2 = 2.555596113204956

$ws.Range("A3").Value2 = "model_8_1_22"
$ws.Range("B3").Value2 = 0.296839247183847
$ws.Range("C3").Value2 = -5.076744665424188
$ws.Range("D3").Value2 = -2.273539485578954
$ws.Range("E3").Value2 = -3.236178887580762
$ws.Range("F3").Value2 = 0.7781917452812195
$ws.Range("G3").Value2 = 2.754521131515503
$ws.Range("H3").Value2 = 2.331639766693115
$ws.Range("I3").Value2 = 2.55551815032959

$ws.Range("A4").Value2 = "model_8_1_20"
$ws.Range("B4").Value2 = 0.3292169385419875
$ws.Range("C4").Value2 = -5.110838622375698
$ws.Range("D4").Value2 = -1.895658530372866
$ws.Range("E4").Value2 = -3.039783100378081
$ws.Range("F4").Value2 = 0.7423592209815979
$ws.Range("G4").Value2 = 2.769975185394287
$ws.Range("H4").Value2 = 2.062487125396729
$ws.Range("I4").Value2 = 2.437040328979492

$ws.Range("A5").Value2 = "model_8_1_17"
$ws.Range("B5").Value2 = 0.3324699184755932
$ws.Range("C5").Value2 = -4.676507709943019
$ws.Range("D5").Value2 = -2.15443543580422
$ws.Range("E5").Value2 = -3.010787281620557
$ws.Range("F5").Value2 = 0.7387591600418091
$ws.Range("G5").Value2 = 2.573097944259644
$ws.Range("H5").Value2 = 2.246805667877197
$ws.Range("I5").Value2 = 2.419548511505127

$ws.Range("A6").Value2 = "model_8_1_16"
$ws.Range("B6").Value2 = 0.3327148842938553
$ws.Range("C6").Value2 = -4.674472222745455
$ws.Range("D6").Value2 = -2.149080756355334
$ws.Range("E6").Value2 = -3.007001415175568
$ws.Range("F6").Value2 = 0.7384880185127258
$ws.Range("G6").Value2 = 2.572175025939941
$ws.Range("H6").Value2 = 2.24299168586731
$ws.Range("I6").Value2 = 2.417264699935913

$ws.Range("A7").Value2 = "model_8_1_18"
$ws.Range("B7").Value2 = 0.3365520873084841
$ws.Range("C7").Value2 = -4.605686881728738
$ws.Range("D7").Value2 = -2.163475883280205
$ws.Range("E7").Value2 = -2.987637870432666
$ws.Range("F7").Value2 = 0.7342414259910583
$ws.Range("G7").Value2 = 2.540995597839355
$ws.Range("H7").Value2 = 2.253244876861572
$ws.Range("I7").Value2 = 2.405583143234253

$ws.Range("A8").Value2 = "model_8_1_23"
$ws.Range("B8").Value2 = 0.3398899891663768
$ws.Range("C8").Value2 = -5.043193315488581
$ws.Range("D8").Value2 = -1.829759389785403
$ws.Range("E8").Value2 = -2.976250263278335
$ws.Range("F8").Value2 = 0.7305472493171692
$ws.Range("G8").Value2 = 2.739312410354614
$ws.Range("H8").Value2 = 2.015548944473267
$ws.Range("I8").Value2 = 2.3987135887146

$ws.Range("A9").Value2 = "model_8_1_24"
$ws.Range("B9").Value2 = 0.3762166622488349
$ws.Range("C9").Value2 = -5.013204840902612
$ws.Range("D9").Value2 = -1.456530259647425
$ws.Range("E9").Value2 = -2.756954933486413
$ws.Range("F9").Value2 = 0.6903443336486816
$ws.Range("G9").Value2 = 2.725718975067139
$ws.Range("H9").Value2 = 1.749709606170654
$ws.Range("I9").Value2 = 2.266421318054199

$ws.Range("A10").Value2 = "model_8_1_19"
$ws.Range("B10").Value2 = 0.3938622197349974
$ws.Range("C10").Value2 = -5.238591322032885
$ws.Range("D10").Value2 = -1.097502855301441
$ws.Range("E10").Value2 = -2.647128048319654
$ws.Range("F10").Value2 = 0.6708158850669861
$ws.Range("G10").Value2 = 2.827884197235107
$ws.Range("H10").Value2 = 1.493985652923584
$ws.Range("I10").Value2 = 2.200167417526245

$ws.Range("A11").Value2 = "model_8_1_15"
$ws.Range("B11").Value2 = 0.4112606094597392
$ws.Range("C11").Value2 = -4.237789433697467
$ws.Range("D11").Value2 = -1.55942423687763
$ws.Range("E11").Value2 = -2.50566490660805
$ws.Range("F11").Value2 = 0.65156090259552
$ws.Range("G11").Value2 = 2.374232053756714
$ws.Range("H11").Value2 = 1.82299792766571
$ws.Range("I11").Value2 = 2.114828109741211

$ws.Range("A12").Value2 = "model_8_1_14"
$ws.Range("B12").Value2 = 0.4877052816509916
$ws.Range("C12").Value2 = -3.528782244898188
$ws.Range("D12").Value2 = -1.060722176129559
$ws.Range("E12").Value2 = -1.94653326459686
$ws.Range("F12").Value2 = 0.5669593214988708
$ws.Range("G12").Value2 = 2.052846670150757
$ws.Range("H12").Value2 = 1.467787981033325
$ws.Range("I12").Value2 = 1.777526259422302

$ws.Range("A13").Value2 = "model_8_1_13"
$ws.Range("B13").Value2 = 0.5123168459875944
$ws.Range("C13").Value2 = -2.959773867237705
$ws.Range("D13").Value2 = -1.195922964855618
$ws.Range("E13").Value2 = -1.795302246636174
$ws.Range("F13").Value2 = 0.5397215485572815
$ws.Range("G13").Value2 = 1.794921636581421
$ws.Range("H13").Value2 = 1.564087271690369
$ws.Range("I13").Value2 = 1.686294555664062

$ws.Range("A14").Value2 = "model_8_1_0"
$ws.Range("B14").Value2 = 0.5238340768688178
$ws.Range("C14").Value2 = 0.03218976656664285
$ws.Range("D14").Value2 = 0.7996597071090313
$ws.Range("E14").Value2 = 0.5036922289990142
$ws.Range("F14").Value2 = 0.5269753336906433
$ws.Range("G14").Value2 = 0.4386976361274719
$ws.Range("H14").Value2 = 0.1426961272954941
$ws.Range("I14").Value2 = 0.2994027137756348

$ws.Range("A15").Value2 = "model_8_1_1"
$ws.Range("B15").Value2 = 0.5322852003042579
$ws.Range("C15").Value2 = 0.1439067621487059
$ws.Range("D15").Value2 = 0.575950183063197
$ws.Range("E15").Value2 = 0.423836488349761
$ws.Range("F15").Value2 = 0.5176224708557129
$ws.Range("G15").Value2 = 0.3880575597286224
$ws.Range("H15").Value2 = 0.3020374178886414
$ws.Range("I15").Value2 = 0.3475765287876129

$ws.Range("A16").Value2 = "model_8_1_2"
$ws.Range("B16").Value2 = 0.5562572988000918
$ws.Range("C16").Value2 = 0.1760896164881085
$ws.Range("D16").Value2 = 0.5602927074744977
$ws.Range("E16").Value2 = 0.427940826398951
$ws.Range("F16").Value2 = 0.4910923838615417
$ws.Range("G16").Value2 = 0.3734694719314575
$ws.Range("H16").Value2 = 0.3131897449493408
$ws.Range("I16").Value2 = 0.3451005220413208

$ws.Range("A17").Value2 = "model_8_1_3"
$ws.Range("B17").Value2 = 0.5947990555705347
$ws.Range("C17").Value2 = 0.1632972725320558
$ws.Range("D17").Value2 = 0.5168765438938188
$ws.Range("E17").Value2 = 0.3987252622558896
$ws.Range("F17").Value2 = 0.4484380185604095
$ws.Range("G17").Value2 = 0.3792680501937866
$ws.Range("H17").Value2 = 0.3441137373447418
$ws.Range("I17").Value2 = 0.3627251088619232

$ws.Range("A18").Value2 = "model_8_1_12"
$ws.Range("B18").Value2 = 0.6007792680896009
$ws.Range("C18").Value2 = -1.938149122355203
$ws.Range("D18").Value2 = -0.8619915416496553
$ws.Range("E18").Value2 = -1.203359331668816
$ws.Range("F18").Value2 = 0.4418197274208069
$ws.Range("G18").Value2 = 1.331830501556396
$ws.Range("H18").Value2 = 1.326238393783569
$ws.Range("I18").Value2 = 1.329198956489563

$ws.Range("A19").Value2 = "model_8_1_11"
$ws.Range("B19").Value2 = 0.6262386150756667
$ws.Range("C19").Value2 = -1.761868791054106
$ws.Range("D19").Value2 = -0.6828014845847421
$ws.Range("E19").Value2 = -1.033675053999164
$ws.Range("F19").Value2 = 0.4136436879634857
$ws.Range("G19").Value2 = 1.251924514770508
$ws.Range("H19").Value2 = 1.198606848716736
$ws.Range("I19").Value2 = 1.226835131645203

$ws.Range("A20").Value2 = "model_8_1_4"
$ws.Range("B20").Value2 = 0.6836157563672909
$ws.Range("C20").Value2 = 0.2606719134860825
$ws.Range("D20").Value2 = -0.428755104622103
$ws.Range("E20").Value2 = -0.08795120517621235
$ws.Range("F20").Value2 = 0.3501441180706024
$ws.Range("G20").Value2 = 0.3351292014122009
$ws.Range("H20").Value2 = 1.017657518386841
$ws.Range("I20").Value2 = 0.6563176512718201

$ws.Range("A21").Value2 = "model_8_1_5"
$ws.Range("B21").Value2 = 0.7306517394767716
$ws.Range("C21").Value2 = 0.06406609034753108
$ws.Range("D21").Value2 = -0.3496078834755589
$ws.Range("E21").Value2 = -0.1221868161700823
$ws.Range("F21").Value2 = 0.2980891764163971
$ws.Range("G21").Value2 = 0.4242484569549561
$ws.Range("H21").Value2 = 0.9612834453582764
$ws.Range("I21").Value2 = 0.6769706606864929

$ws.Range("A22").Value2 = "model_8_1_7"
$ws.Range("B22").Value2 = 0.7492582075702461
$ws.Range("C22").Value2 = 0.06788172426717554
$ws.Range("D22").Value2 = -0.2295331349895091
$ws.Range("E22").Value2 = -0.05395136376800269
$ws.Range("F22").Value2 = 0.277497261762619
$ws.Range("G22").Value2 = 0.4225188791751862
$ws.Range("H22").Value2 = 0.8757579922676086
$ws.Range("I22").Value2 = 0.6358069181442261

$ws.Range("A23").Value2 = "model_8_1_6"
$ws.Range("B23").Value2 = 0.7494843026606581
$ws.Range("C23").Value2 = 0.08176120355294425
$ws.Range("D23").Value2 = -0.1799332347207638
$ws.Range("E23").Value2 = -0.0208713164356471
$ws.Range("F23").Value2 = 0.2772470414638519
$ws.Range("G23").Value2 = 0.4162274599075317
$ws.Range("H23").Value2 = 0.8404295444488525
$ws.Range("I23").Value2 = 0.6158509850502014

$ws.Range("A24").Value2 = "model_8_1_10"
$ws.Range("B24").Value2 = 0.7608485686723295
$ws.Range("C24").Value2 = -0.2332293720618914
$ws.Range("D24").Value2 = -0.1652522724791425
$ws.Range("E24").Value2 = -0.138019816863439
$ws.Range("F24").Value2 = 0.2646701633930206
$ws.Range("G24").Value2 = 0.5590091943740845
$ws.Range("H24").Value2 = 0.8299726843833923
$ws.Range("I24").Value2 = 0.6865220665931702

$ws.Range("A25").Value2 = "model_8_1_9"
$ws.Range("B25").Value2 = 0.7794848471683788
$ws.Range("C25").Value2 = -0.0793253862502914
$ws.Range("D25").Value2 = -0.05429924085652704
$ws.Range("E25").Value2 = -0.01514890376824951
$ws.Range("F25").Value2 = 0.2440452873706818
$ws.Range("G25").Value2 = 0.4892462193965912
$ws.Range("H25").Value2 = 0.7509443759918213
$ws.Range("I25").Value2 = 0.6123989224433899

$ws.Range("A26").Value2 = "model_8_1_8"
$ws.Range("B26").Value2 = 0.7870004994014785
$ws.Range("C26").Value2 = 0.02008300777991379
$ws.Range("D26").Value2 = 0.002122587302799128
$ws.Range("E26").Value2 = 0.05574581868173767
$ws.Range("F26").Value2 = 0.2357276827096939
$ws.Range("G26").Value2 = 0.4441854655742645
$ws.Range("H26").Value2 = 0.7107568979263306
$ws.Range("I26").Value2 = 0.5696309804916382
